# Apply updated NATMI computation values (following Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$data = @{
    2  = @{ E=3; G=8.379101; H=25.137303; I=0.2232365200207407; J=0.2232365200207407; K=3; M=33.24999533333333; N=99.74998599999999; O=0.5673360890306117; P=0.5673360890306117; Q=278.6050691475286; R=2507.445622327758; S=0.1266501341973709; T=0.1266501341973709 }
    3  = @{ E=3; G=8.379101; H=25.137303; I=0.2232365200207407; J=0.2232365200207407; K=3; M=23.30243966666667; N=69.907319; O=0.3976035140102714; P=0.3976035140102714; Q=195.2534955134064; R=1757.281459620657; S=0.0887596248156708; T=0.08875962481567079 }
    4  = @{ E=3; G=8.379101; H=25.137303; I=0.2232365200207407; J=0.2232365200207407; K=3; M=2.054792666666667; N=6.164378; O=0.03506039695911681; P=0.03506039695911681; Q=17.21731528805934; R=154.955837592534; S=0.007826761007698996; T=0.007826761007698994 }
    5  = @{ E=3; G=9.350178333333334; H=28.050535; I=0.2491080215773342; J=0.2491080215773342; K=3; M=33.24999533333333; N=99.74998599999999; O=0.5673360890306117; P=0.5673360890306117; Q=310.8933859491678; R=2798.04047354251; S=0.141327970707838; T=0.141327970707838 }
    6  = @{ E=3; G=9.350178333333334; H=28.050535; I=0.2491080215773342; J=0.2491080215773342; K=3; M=23.30243966666667; N=69.907319; O=0.3976035140102714; P=0.3976035140102714; Q=217.8819664850739; R=1960.937698365665; S=0.0990462247472946; T=0.09904622474729458 }
    7  = @{ E=3; G=9.350178333333334; H=28.050535; I=0.2491080215773342; J=0.2491080215773342; K=3; M=2.054792666666667; N=6.164378; O=0.03506039695911681; P=0.03506039695911681; Q=19.21267787135889; R=172.91410084223; S=0.008733826122201574; T=0.00873382612220157 }
    8  = @{ E=3; G=19.80535433333333; H=59.416063; I=0.5276554584019252; J=0.5276554584019252; K=3; M=33.24999533333333; N=99.74998599999999; O=0.5673360890306117; P=0.5673360890306117; Q=658.5279391583464; R=5926.751452425117; S=0.2993579841254029; T=0.2993579841254029 }
    9  = @{ E=3; G=19.80535433333333; H=59.416063; I=0.5276554584019252; J=0.5276554584019252; K=3; M=23.30243966666667; N=69.907319; O=0.3976035140102714; P=0.3976035140102714; Q=461.5130744294553; R=4153.617669865097; S=0.209797664447306; T=0.209797664447306 }
    10 = @{ E=3; G=19.80535433333333; H=59.416063; I=0.5276554584019252; J=0.5276554584019252; K=3; M=2.054792666666667; N=6.164378; O=0.03506039695911681; P=0.03506039695911681; Q=40.69589684486822; R=366.263071603814; S=0.01849980982921625; T=0.01849980982921624 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
